# Add a new row (row 9) of department data to the "部门情况202401" sheet
# (the first worksheet), which currently spans A1:O8. The new row holds
# the "数字普惠部" department figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 9
$rng = $ws.Range("A$($newRow):O$($newRow)")

# Force the cells to be stored as text so that numeric-looking values such
# as "0.00" keep their literal textual representation instead of being
# coerced into numbers.
$rng.NumberFormat = "@"

$values = @(
    "数字普惠部",
    "0.00",
    "0.00",
    "0.00",
    "0.00",
    "95076.05",
    "6240.00",
    "18.66",
    "22305.70",
    "0.00",
    "0.00",
    "0.00",
    "0.00",
    "0.00",
    "0.00"
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")

for ($i = 0; $i -lt $columns.Length; $i++) {
    $cellAddr = "$($columns[$i])$($newRow)"
    $ws.Range($cellAddr).Value = $values[$i]
}
